$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.903.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.770.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.11%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.27%  '

$ws.Range("E7").Value = '  +3.44%  '

$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.793.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.43%  '

$ws.Range("E10").Value = '  +3.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.399'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.27%  '

$ws.Range("E13").Value = '  +1.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.264.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.849.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.63%  '

$ws.Range("E17").Value = '  +7.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.784.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '368.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.560'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.73%  '

$ws.Range("E24").Value = '  +0.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.83%  '

$ws.Range("E26").Value = '  +6.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0968'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +16.51%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '

$ws.Range("E37").Value = '  +7.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.92%  '

$ws.Range("E39").Value = '  +3.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '344.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.28'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0613'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.654'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.68%  '

$ws.Range("E48").Value = '  +2.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.15%  '

$ws.Range("E50").Value = '  +2.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.177.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
